$d = $word.ActiveDocument

$d.Content.Find.Execute("27+44=71", $true, $true, $false, $false, $false, $true, 1, $false, "35+17=52", 2) | Out-Null
$d.Content.Find.Execute("37+2=39", $true, $true, $false, $false, $false, $true, 1, $false, "95-5=90", 2) | Out-Null
$d.Content.Find.Execute("64-3=61", $true, $true, $false, $false, $false, $true, 1, $false, "48-26=22", 2) | Out-Null
$d.Content.Find.Execute("5+11=16", $true, $true, $false, $false, $false, $true, 1, $false, "89-75=14", 2) | Out-Null
$d.Content.Find.Execute("45-35=10", $true, $true, $false, $false, $false, $true, 1, $false, "16-1=15", 2) | Out-Null
$d.Content.Find.Execute("46+51=97", $true, $true, $false, $false, $false, $true, 1, $false, "75-41=34", 2) | Out-Null
$d.Content.Find.Execute("43+20=63", $true, $true, $false, $false, $false, $true, 1, $false, "15+56=71", 2) | Out-Null
$d.Content.Find.Execute("91-42=49", $true, $true, $false, $false, $false, $true, 1, $false, "8+73=81", 2) | Out-Null
$d.Content.Find.Execute("52+33=85", $true, $true, $false, $false, $false, $true, 1, $false, "51-23=28", 2) | Out-Null
$d.Content.Find.Execute("59+35=94", $true, $true, $false, $false, $false, $true, 1, $false, "57-22=35", 2) | Out-Null
$d.Content.Find.Execute("25+13=38", $true, $true, $false, $false, $false, $true, 1, $false, "7+71=78", 2) | Out-Null
$d.Content.Find.Execute("69+20=89", $true, $true, $false, $false, $false, $true, 1, $false, "37-29=8", 2) | Out-Null
$d.Content.Find.Execute("18-15=3", $true, $true, $false, $false, $false, $true, 1, $false, "10+57=67", 2) | Out-Null
$d.Content.Find.Execute("21-19=2", $true, $true, $false, $false, $false, $true, 1, $false, "50-38=12", 2) | Out-Null
$d.Content.Find.Execute("74-12=62", $true, $true, $false, $false, $false, $true, 1, $false, "7+83=90", 2) | Out-Null
$d.Content.Find.Execute("76-31=45", $true, $true, $false, $false, $false, $true, 1, $false, "32+38=70", 2) | Out-Null
$d.Content.Find.Execute("13+85=98", $true, $true, $false, $false, $false, $true, 1, $false, "62+16=78", 2) | Out-Null
$d.Content.Find.Execute("27+51=78", $true, $true, $false, $false, $false, $true, 1, $false, "12-2=10", 2) | Out-Null
$d.Content.Find.Execute("5+24=29", $true, $true, $false, $false, $false, $true, 1, $false, "57-50=7", 2) | Out-Null
$d.Content.Find.Execute("1+45=46", $true, $true, $false, $false, $false, $true, 1, $false, "29+62=91", 2) | Out-Null
$d.Content.Find.Execute("24+1=25", $true, $true, $false, $false, $false, $true, 1, $false, "84-13=71", 2) | Out-Null
$d.Content.Find.Execute("16+82=98", $true, $true, $false, $false, $false, $true, 1, $false, "25+17=42", 2) | Out-Null
$d.Content.Find.Execute("6+42=48", $true, $true, $false, $false, $false, $true, 1, $false, "21-0=21", 2) | Out-Null
$d.Content.Find.Execute("83-38=45", $true, $true, $false, $false, $false, $true, 1, $false, "76-75=1", 2) | Out-Null
$d.Content.Find.Execute("26-14=12", $true, $true, $false, $false, $false, $true, 1, $false, "17+61=78", 2) | Out-Null
$d.Content.Find.Execute("31-23=8", $true, $true, $false, $false, $false, $true, 1, $false, "32-14=18", 2) | Out-Null
$d.Content.Find.Execute("42+43=85", $true, $true, $false, $false, $false, $true, 1, $false, "32-24=8", 2) | Out-Null
$d.Content.Find.Execute("76-39=37", $true, $true, $false, $false, $false, $true, 1, $false, "9-2=7", 2) | Out-Null
$d.Content.Find.Execute("25+64=89", $true, $true, $false, $false, $false, $true, 1, $false, "57+9=66", 2) | Out-Null
$d.Content.Find.Execute("69-64=5", $true, $true, $false, $false, $false, $true, 1, $false, "97-73=24", 2) | Out-Null
$d.Content.Find.Execute("61+20=81", $true, $true, $false, $false, $false, $true, 1, $false, "30+51=81", 2) | Out-Null
$d.Content.Find.Execute("42-34=8", $true, $true, $false, $false, $false, $true, 1, $false, "65-12=53", 2) | Out-Null
$d.Content.Find.Execute("14+44=58", $true, $true, $false, $false, $false, $true, 1, $false, "99-98=1", 2) | Out-Null
$d.Content.Find.Execute("65-26=39", $true, $true, $false, $false, $false, $true, 1, $false, "90-56=34", 2) | Out-Null
$d.Content.Find.Execute("60+9=69", $true, $true, $false, $false, $false, $true, 1, $false, "19+23=42", 2) | Out-Null
$d.Content.Find.Execute("15+46=61", $true, $true, $false, $false, $false, $true, 1, $false, "7+31=38", 2) | Out-Null
$d.Content.Find.Execute("4+65=69", $true, $true, $false, $false, $false, $true, 1, $false, "1+64=65", 2) | Out-Null
$d.Content.Find.Execute("59+7=66", $true, $true, $false, $false, $false, $true, 1, $false, "6+79=85", 2) | Out-Null
$d.Content.Find.Execute("85-80=5", $true, $true, $false, $false, $false, $true, 1, $false, "78-73=5", 2) | Out-Null
$d.Content.Find.Execute("9+82=91", $true, $true, $false, $false, $false, $true, 1, $false, "69+17=86", 2) | Out-Null
$d.Content.Find.Execute("54+26=80", $true, $true, $false, $false, $false, $true, 1, $false, "47+50=97", 2) | Out-Null
$d.Content.Find.Execute("22+10=32", $true, $true, $false, $false, $false, $true, 1, $false, "98-96=2", 2) | Out-Null
$d.Content.Find.Execute("77-54=23", $true, $true, $false, $false, $false, $true, 1, $false, "6+31=37", 2) | Out-Null
$d.Content.Find.Execute("40-21=19", $true, $true, $false, $false, $false, $true, 1, $false, "38+32=70", 2) | Out-Null
$d.Content.Find.Execute("1+90=91", $true, $true, $false, $false, $false, $true, 1, $false, "82-35=47", 2) | Out-Null
$d.Content.Find.Execute("58-45=13", $true, $true, $false, $false, $false, $true, 1, $false, "21+58=79", 2) | Out-Null
$d.Content.Find.Execute("77-25=52", $true, $true, $false, $false, $false, $true, 1, $false, "50-10=40", 2) | Out-Null
$d.Content.Find.Execute("82-81=1", $true, $true, $false, $false, $false, $true, 1, $false, "51+24=75", 2) | Out-Null
$d.Content.Find.Execute("73-48=25", $true, $true, $false, $false, $false, $true, 1, $false, "54-19=35", 2) | Out-Null
$d.Content.Find.Execute("85-62=23", $true, $true, $false, $false, $false, $true, 1, $false, "54+8=62", 2) | Out-Null
$d.Content.Find.Execute("41-26=15", $true, $true, $false, $false, $false, $true, 1, $false, "81-31=50", 2) | Out-Null
$d.Content.Find.Execute("1+53=54", $true, $true, $false, $false, $false, $true, 1, $false, "93-6=87", 2) | Out-Null
$d.Content.Find.Execute("6+84=90", $true, $true, $false, $false, $false, $true, 1, $false, "12+44=56", 2) | Out-Null
$d.Content.Find.Execute("47+38=85", $true, $true, $false, $false, $false, $true, 1, $false, "36+55=91", 2) | Out-Null
$d.Content.Find.Execute("20+29=49", $true, $true, $false, $false, $false, $true, 1, $false, "75-34=41", 2) | Out-Null
$d.Content.Find.Execute("43+34=77", $true, $true, $false, $false, $false, $true, 1, $false, "79-4=75", 2) | Out-Null
$d.Content.Find.Execute("96-79=17", $true, $true, $false, $false, $false, $true, 1, $false, "6+52=58", 2) | Out-Null
$d.Content.Find.Execute("97-50=47", $true, $true, $false, $false, $false, $true, 1, $false, "61+32=93", 2) | Out-Null
$d.Content.Find.Execute("14+27=41", $true, $true, $false, $false, $false, $true, 1, $false, "58+8=66", 2) | Out-Null
$d.Content.Find.Execute("18+71=89", $true, $true, $false, $false, $false, $true, 1, $false, "94-7=87", 2) | Out-Null
$d.Content.Find.Execute("28-3=25", $true, $true, $false, $false, $false, $true, 1, $false, "64-53=11", 2) | Out-Null
$d.Content.Find.Execute("86-71=15", $true, $true, $false, $false, $false, $true, 1, $false, "16-15=1", 2) | Out-Null
$d.Content.Find.Execute("9+20=29", $true, $true, $false, $false, $false, $true, 1, $false, "54-44=10", 2) | Out-Null
$d.Content.Find.Execute("33+55=88", $true, $true, $false, $false, $false, $true, 1, $false, "68-59=9", 2) | Out-Null
$d.Content.Find.Execute("31+14=45", $true, $true, $false, $false, $false, $true, 1, $false, "34+53=87", 2) | Out-Null
$d.Content.Find.Execute("16+56=72", $true, $true, $false, $false, $false, $true, 1, $false, "35+17=52", 2) | Out-Null
$d.Content.Find.Execute("28-11=17", $true, $true, $false, $false, $false, $true, 1, $false, "79-7=72", 2) | Out-Null
$d.Content.Find.Execute("10+82=92", $true, $true, $false, $false, $false, $true, 1, $false, "54-32=22", 2) | Out-Null
$d.Content.Find.Execute("89-65=24", $true, $true, $false, $false, $false, $true, 1, $false, "38-16=22", 2) | Out-Null
$d.Content.Find.Execute("12+83=95", $true, $true, $false, $false, $false, $true, 1, $false, "48-1=47", 2) | Out-Null
$d.Content.Find.Execute("26+46=72", $true, $true, $false, $false, $false, $true, 1, $false, "90-33=57", 2) | Out-Null
$d.Content.Find.Execute("3+48=51", $true, $true, $false, $false, $false, $true, 1, $false, "17-2=15", 2) | Out-Null
$d.Content.Find.Execute("9+89=98", $true, $true, $false, $false, $false, $true, 1, $false, "17+34=51", 2) | Out-Null
$d.Content.Find.Execute("77-3=74", $true, $true, $false, $false, $false, $true, 1, $false, "51-47=4", 2) | Out-Null
$d.Content.Find.Execute("89+10=99", $true, $true, $false, $false, $false, $true, 1, $false, "82-40=42", 2) | Out-Null
$d.Content.Find.Execute("81-50=31", $true, $true, $false, $false, $false, $true, 1, $false, "56-4=52", 2) | Out-Null
$d.Content.Find.Execute("65+10=75", $true, $true, $false, $false, $false, $true, 1, $false, "99-38=61", 2) | Out-Null
$d.Content.Find.Execute("59-12=47", $true, $true, $false, $false, $false, $true, 1, $false, "56+33=89", 2) | Out-Null
$d.Content.Find.Execute("20-14=6", $true, $true, $false, $false, $false, $true, 1, $false, "55-47=8", 2) | Out-Null
$d.Content.Find.Execute("53-5=48", $true, $true, $false, $false, $false, $true, 1, $false, "44+25=69", 2) | Out-Null
$d.Content.Find.Execute("69-0=69", $true, $true, $false, $false, $false, $true, 1, $false, "48-4=44", 2) | Out-Null
$d.Content.Find.Execute("15+9=24", $true, $true, $false, $false, $false, $true, 1, $false, "57+28=85", 2) | Out-Null
$d.Content.Find.Execute("52+40=92", $true, $true, $false, $false, $false, $true, 1, $false, "63+32=95", 2) | Out-Null
$d.Content.Find.Execute("52-45=7", $true, $true, $false, $false, $false, $true, 1, $false, "94-39=55", 2) | Out-Null
$d.Content.Find.Execute("89-57=32", $true, $true, $false, $false, $false, $true, 1, $false, "51+27=78", 2) | Out-Null
$d.Content.Find.Execute("9-5=4", $true, $true, $false, $false, $false, $true, 1, $false, "0+72=72", 2) | Out-Null
$d.Content.Find.Execute("37+20=57", $true, $true, $false, $false, $false, $true, 1, $false, "85-67=18", 2) | Out-Null
$d.Content.Find.Execute("8+57=65", $true, $true, $false, $false, $false, $true, 1, $false, "43+4=47", 2) | Out-Null
$d.Content.Find.Execute("52+22=74", $true, $true, $false, $false, $false, $true, 1, $false, "98-73=25", 2) | Out-Null
$d.Content.Find.Execute("80-45=35", $true, $true, $false, $false, $false, $true, 1, $false, "6+38=44", 2) | Out-Null
$d.Content.Find.Execute("53+18=71", $true, $true, $false, $false, $false, $true, 1, $false, "2+86=88", 2) | Out-Null
$d.Content.Find.Execute("92-22=70", $true, $true, $false, $false, $false, $true, 1, $false, "16+23=39", 2) | Out-Null
$d.Content.Find.Execute("93-37=56", $true, $true, $false, $false, $false, $true, 1, $false, "85-71=14", 2) | Out-Null
$d.Content.Find.Execute("63-32=31", $true, $true, $false, $false, $false, $true, 1, $false, "78-16=62", 2) | Out-Null
$d.Content.Find.Execute("15+60=75", $true, $true, $false, $false, $false, $true, 1, $false, "64-12=52", 2) | Out-Null
$d.Content.Find.Execute("57-2=55", $true, $true, $false, $false, $false, $true, 1, $false, "79-45=34", 2) | Out-Null
$d.Content.Find.Execute("39+56=95", $true, $true, $false, $false, $false, $true, 1, $false, "70-58=12", 2) | Out-Null
$d.Content.Find.Execute("35-16=19", $true, $true, $false, $false, $false, $true, 1, $false, "71-23=48", 2) | Out-Null
$d.Content.Find.Execute("56+16=72", $true, $true, $false, $false, $false, $true, 1, $false, "13+84=97", 2) | Out-Null
$d.Content.Find.Execute("60-41=19", $true, $true, $false, $false, $false, $true, 1, $false, "74-21=53", 2) | Out-Null
